$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 9) with the next day's data, mirroring row 8's values
$ws.Range("A9").Value = "'2025-08-24"
$ws.Range("A9").Style = "Normal"
$ws.Range("B9").Value = 58.5099983215332
$ws.Range("C9").Value = 680.2999877929688
$ws.Range("D9").Value = 319.1000061035156
